$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("AddOpportunity")
$ws.Range("R4").Value = "Public Equity"
$ws.Range("R5").Value = "Public Equity"
$ws.Range("B3").Value = "Owl Wire and Cable LLC"
$ws.Range("A3").Value = "International Wire Group, Inc."
$ws.Range("A4").Value = "3i Group Plc"
$ws.Range("B4").Value = "IRISNDT Corp."
$ws.Range("A5").Value = "Mirait Holdings Corporation"
$ws.Range("B5").Value = "Seibu Construction Co., Ltd."
$ws.Columns("A:B").AutoFit()
$ws.Range("B12").Select()
